$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "cryptos" price list (GitHub Actions nightly scrape).
# Only the Price (D) / Volume(1h) (E) columns move for most rows; three
# coin pairs additionally swapped rank position (rows 10/11, 17/18) and
# Decentraland (row 51) was replaced by NEARProtocol in the ranking.
#
# Price values that look numeric (e.g. "239.69") are written with a
# leading apostrophe so Excel stores them as literal text instead of
# re-parsing them as a Double (this data uses "." as a thousands marker
# for some rows, e.g. "30.633.50", so every Price cell must stay text).
# The immediate Style = "Normal" afterwards clears the resulting
# quote-prefix marker so the cell format matches the untouched cells.

$ws.Range("D2").Value = "30.633.50"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "1.894.50"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'239.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.34%  "
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").Value = "'0.4919"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.10%  "
$ws.Range("D8").Value = "'0.2947"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.52%  "
$ws.Range("E9").Value = "  +1.40%  "
$ws.Range("B10").Value = "WrappedEther"
$ws.Range("C10").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D10").Value = "1.927.19"
$ws.Range("E10").Value = "  +2.56%  "
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").Value = "'17.16"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.73%  "
$ws.Range("D12").Value = "'0.07358"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.03%  "
$ws.Range("D13").Value = "'5.159"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.66%  "
$ws.Range("D14").Value = "'88.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.04%  "
$ws.Range("D15").Value = "'0.6716"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.67%  "
$ws.Range("D16").Value = "30.575.92"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").Value = "'13.49"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.36%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.000007889"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.14%  "
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").Value = "2.106.49"
$ws.Range("E20").Value = "  -0.63%  "
$ws.Range("B21").Value = "BinanceUSD"
$ws.Range("C21").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D21").Value = "'1.003"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'5.283"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +12.01%  "
$ws.Range("D23").Value = "'191.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.51%  "
$ws.Range("E24").Value = "  +3.23%  "
$ws.Range("D25").Value = "'9.534"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.95%  "
$ws.Range("D26").Value = "'161.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.93%  "
$ws.Range("D27").Value = "'18.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.73%  "
$ws.Range("D28").Value = "'1.949"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.96%  "
$ws.Range("D29").Value = "'1.477"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.69%  "
$ws.Range("D30").Value = "'4.448"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.03%  "
$ws.Range("D31").Value = "'0.09213"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.43%  "
$ws.Range("D32").Value = "'4.154"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.21%  "
$ws.Range("D33").Value = "'0.05248"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.11%  "
$ws.Range("D34").Value = "'0.7437"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.63%  "
$ws.Range("D35").Value = "'1.108"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.11%  "
$ws.Range("E36").Value = "  +0.98%  "
$ws.Range("D37").Value = "'0.01839"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.14%  "
$ws.Range("D38").Value = "'2.697"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.82%  "
$ws.Range("D39").Value = "'0.9254"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.59%  "
$ws.Range("D40").Value = "'2.058"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.13%  "
$ws.Range("D41").Value = "'0.4424"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.98%  "
$ws.Range("D42").Value = "'5.964"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.05%  "
$ws.Range("D43").Value = "'106.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.48%  "
$ws.Range("D44").Value = "'71.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +25.52%  "
$ws.Range("D45").Value = "'0.9942"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("E46").Value = "  +4.17%  "
$ws.Range("D47").Value = "'7.595"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.92%  "
$ws.Range("D48").Value = "'9.063"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.17%  "
$ws.Range("D49").Value = "'35.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.45%  "
$ws.Range("D50").Value = "'0.05826"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "'1.430"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.55%  "
